# Add a new "product_type" column (F) to the checklist sheet, classifying
# each product as a rulebook, replay, or supplement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("F1").Value = "product_type"

# Row data: year / japanese / english / publisher / image already exist in A:E.
# Populate the new product_type column per row.
$ws.Range("F2").Value = "rulebook"
$ws.Range("F3").Value = "replay"
$ws.Range("F4").Value = "replay"
$ws.Range("F5").Value = "replay"
$ws.Range("F6").Value = "rulebook"
$ws.Range("F7").Value = "supplement"
$ws.Range("F8").Value = "replay"

# Match the author's final selection state (cell F9, just below the new
# column's last populated row).
$ws.Range("F9").Select()
